{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" copyright/footer line, and the blank paragraph\n// that separates them from the preceding \"LOT2041: ...\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetFollow = \"LOT2041: Engenharia Bioqu\u00edmica (Requisito fraco)\";\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetFollow) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The blank paragraph, the \"Ver no Jupiter...\" paragraph, and the\n  // \"\u00a9 2020...\" paragraph immediately follow the requisito paragraph.\n  const blank = items[anchorIndex + 1];\n  const jupiter = items[anchorIndex + 2];\n  const copyright = items[anchorIndex + 3];\n\n  if (blank && blank.text === \"\" && jupiter && jupiter.text === jupiterText && copyright && copyright.text === copyrightText) {\n    // Delete from last to first so earlier indices stay valid.\n    copyright.delete();\n    jupiter.delete();\n    blank.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" copyright/footer line, and the blank paragraph\n# that separates them from the preceding \"LOT2041: ...\" requirement line.\n$d = $word.ActiveDocument\n\n$anchorText = \"LOT2041: Engenharia Bioqu\" + [char]0x00ED + \"mica (Requisito fraco)\"\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$n = $d.Paragraphs.Count\nfor ($i = $n; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $blank = $d.Paragraphs.Item($i + 1)\n        $jupiter = $d.Paragraphs.Item($i + 2)\n        $copyright = $d.Paragraphs.Item($i + 3)\n\n        $blankText = $blank.Range.Text.TrimEnd([char]13, [char]7)\n        $jupiterFound = $jupiter.Range.Text.TrimEnd([char]13, [char]7)\n        $copyrightFound = $copyright.Range.Text.TrimEnd([char]13, [char]7)\n\n        if ($blankText -eq \"\" -and $jupiterFound -eq $jupiterText -and $copyrightFound -eq $copyrightText) {\n            # Delete from last to first so earlier indices stay valid.\n            $copyright.Range.Delete()\n            $jupiter.Range.Delete()\n            $blank.Range.Delete()\n        }\n        break\n    }\n}\n"}
